$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (shifts old row2->row3, row3->row4)
$ws.Rows.Item(2).Insert()

# New row 2: Minciencias postdoc info
$ws.Range("A2").Value = "Investigadora Principal"
$ws.Range("B2").Value = "Desde 2023"
$ws.Range("C2").Value = "Asociación Red de Mujeres Víctimas y Profesionales"
$ws.Range("D2").Value = "Bogotá, Colombia"
$ws.Range("E2").Value = "La necesidad de generar procesos de reparación social a las mujeres víctimas y sobrevivientes de violencias sexuales en el marco del conflicto armado desde el quehacer periodístico. Diversas propuestas de tratamiento según contextos"

# Row 4 (was row3 before insert): clean up the "why" text to a single line
$ws.Range("E4").Value = "El quehacer periodístico en Colombia y sus aporte en los procesos de memoria histórica  en los casos de violencia sexual contra mujeres en Colombia"

$ws.Range("E10").Select()
